$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = 0.1662817551963048
$ws.Cells.Item(2, 3).Value = 0.605080831408776
$ws.Cells.Item(2, 10).Value = 0.0138568129330254
$ws.Cells.Item(2, 16).Value = 0.1120092378752887
$ws.Cells.Item(2, 19).Value = 0.1027713625866051
$ws.Cells.Item(3, 2).Value = 0.01279707495429616
$ws.Cells.Item(3, 3).Value = 0.03290676416819013
$ws.Cells.Item(3, 10).Value = 0.03290676416819013
$ws.Cells.Item(3, 16).Value = 0.7568555758683729
$ws.Cells.Item(3, 19).Value = 0.1645338208409506
$ws.Cells.Item(4, 10).Value = 0.07482993197278912
$ws.Cells.Item(4, 15).Value = 0.006802721088435374
$ws.Cells.Item(4, 16).Value = 0.673469387755102
$ws.Cells.Item(4, 19).Value = 0.2448979591836735
$ws.Cells.Item(6, 2).Value = 0.06230529595015576
$ws.Cells.Item(6, 4).Value = 0.009345794392523364
$ws.Cells.Item(6, 5).Value = 0.003115264797507788
$ws.Cells.Item(6, 6).Value = 0.06697819314641744
$ws.Cells.Item(6, 10).Value = 0.2414330218068536
$ws.Cells.Item(6, 15).Value = 0.01246105919003115
$ws.Cells.Item(6, 17).Value = 0.1526479750778816
$ws.Cells.Item(6, 18).Value = 0.06542056074766354
$ws.Cells.Item(6, 19).Value = 0.3862928348909657
$ws.Cells.Item(7, 2).Value = 0.1252371916508539
$ws.Cells.Item(7, 4).Value = 0.02466793168880456
$ws.Cells.Item(7, 5).Value = 0.00189753320683112
$ws.Cells.Item(7, 6).Value = 0.05123339658444023
$ws.Cells.Item(7, 10).Value = 0.1195445920303605
$ws.Cells.Item(7, 15).Value = 0.03225806451612903
$ws.Cells.Item(7, 17).Value = 0.1859582542694497
$ws.Cells.Item(7, 18).Value = 0.07590132827324478
$ws.Cells.Item(7, 19).Value = 0.3833017077798861
$ws.Cells.Item(8, 2).Value = 0.1063492063492063
$ws.Cells.Item(8, 4).Value = 0.02301587301587302
$ws.Cells.Item(8, 6).Value = 0.0761904761904762
$ws.Cells.Item(8, 10).Value = 0.1126984126984127
$ws.Cells.Item(8, 15).Value = 0.02301587301587302
$ws.Cells.Item(8, 17).Value = 0.1841269841269841
$ws.Cells.Item(8, 18).Value = 0.1031746031746032
$ws.Cells.Item(8, 19).Value = 0.3714285714285714
$ws.Cells.Item(9, 2).Value = 0.1132743362831858
$ws.Cells.Item(9, 4).Value = 0.0247787610619469
$ws.Cells.Item(9, 5).Value = 0.001769911504424779
$ws.Cells.Item(9, 6).Value = 0.06902654867256637
$ws.Cells.Item(9, 10).Value = 0.08849557522123894
$ws.Cells.Item(9, 15).Value = 0.0247787610619469
$ws.Cells.Item(9, 17).Value = 0.1592920353982301
$ws.Cells.Item(9, 18).Value = 0.09911504424778761
$ws.Cells.Item(9, 19).Value = 0.4194690265486726
$ws.Cells.Item(10, 2).Value = 0.1092413793103448
$ws.Cells.Item(10, 4).Value = 0.02455172413793104
$ws.Cells.Item(10, 5).Value = 0.001931034482758621
$ws.Cells.Item(10, 6).Value = 0.07117241379310345
$ws.Cells.Item(10, 10).Value = 0.1097931034482759
$ws.Cells.Item(10, 15).Value = 0.02179310344827586
$ws.Cells.Item(10, 17).Value = 0.2107586206896552
$ws.Cells.Item(10, 18).Value = 0.08744827586206896
$ws.Cells.Item(10, 19).Value = 0.3633103448275862
$ws.Cells.Item(11, 7).Value = 0.1371629542790152
$ws.Cells.Item(11, 10).Value = 0.1148886283704572
$ws.Cells.Item(11, 11).Value = 0.2086752637749121
$ws.Cells.Item(11, 12).Value = 0.52989449003517
$ws.Cells.Item(11, 19).Value = 0.009378663540445486
$ws.Cells.Item(12, 7).Value = 0.7076271186440678
$ws.Cells.Item(12, 10).Value = 0.2097457627118644
$ws.Cells.Item(12, 11).Value = 0.00211864406779661
$ws.Cells.Item(12, 12).Value = 0.02754237288135593
$ws.Cells.Item(12, 19).Value = 0.05296610169491525
$ws.Cells.Item(13, 7).Value = 0.6870229007633588
$ws.Cells.Item(13, 10).Value = 0.2595419847328244
$ws.Cells.Item(13, 19).Value = 0.05343511450381679
$ws.Cells.Item(15, 6).Value = 0.028328611898017
$ws.Cells.Item(15, 8).Value = 0.1628895184135977
$ws.Cells.Item(15, 9).Value = 0.07082152974504249
$ws.Cells.Item(15, 10).Value = 0.3314447592067989
$ws.Cells.Item(15, 11).Value = 0.06657223796033994
$ws.Cells.Item(15, 13).Value = 0.009915014164305949
$ws.Cells.Item(15, 15).Value = 0.06232294617563739
$ws.Cells.Item(15, 19).Value = 0.2677053824362606
$ws.Cells.Item(16, 6).Value = 0.01342281879194631
$ws.Cells.Item(16, 8).Value = 0.1543624161073825
$ws.Cells.Item(16, 9).Value = 0.07718120805369127
$ws.Cells.Item(16, 10).Value = 0.4211409395973154
$ws.Cells.Item(16, 11).Value = 0.09731543624161074
$ws.Cells.Item(16, 13).Value = 0.02516778523489933
$ws.Cells.Item(16, 15).Value = 0.06040268456375839
$ws.Cells.Item(16, 19).Value = 0.151006711409396
$ws.Cells.Item(17, 6).Value = 0.01803921568627451
$ws.Cells.Item(17, 8).Value = 0.1749019607843137
$ws.Cells.Item(17, 9).Value = 0.08941176470588236
$ws.Cells.Item(17, 10).Value = 0.4219607843137255
$ws.Cells.Item(17, 11).Value = 0.09254901960784313
$ws.Cells.Item(17, 13).Value = 0.02274509803921568
$ws.Cells.Item(17, 14).Value = 0.001568627450980392
$ws.Cells.Item(17, 15).Value = 0.0603921568627451
$ws.Cells.Item(17, 19).Value = 0.1184313725490196
$ws.Cells.Item(18, 6).Value = 0.02581755593803787
$ws.Cells.Item(18, 8).Value = 0.197934595524957
$ws.Cells.Item(18, 9).Value = 0.07228915662650602
$ws.Cells.Item(18, 10).Value = 0.4182444061962134
$ws.Cells.Item(18, 11).Value = 0.09122203098106713
$ws.Cells.Item(18, 13).Value = 0.01549053356282272
$ws.Cells.Item(18, 15).Value = 0.06540447504302926
$ws.Cells.Item(18, 19).Value = 0.1135972461273666
$ws.Cells.Item(19, 6).Value = 0.01376791233492554
$ws.Cells.Item(19, 8).Value = 0.2020230401798258
$ws.Cells.Item(19, 9).Value = 0.08738409665636415
$ws.Cells.Item(19, 10).Value = 0.3770722112953077
$ws.Cells.Item(19, 11).Value = 0.1093003652711436
$ws.Cells.Item(19, 13).Value = 0.02107333520651869
$ws.Cells.Item(19, 14).Value = 0.0008429334082607474
$ws.Cells.Item(19, 15).Value = 0.07867378477100309
$ws.Cells.Item(19, 19).Value = 0.1098623208766508
